# gsc-export-old/HTTPS.xlsx update ("updated legacy GSC export data")
#
# The legacy GSC export dropped its oldest day (2025-10-19) from the
# rolling date window on the "Chart" sheet: the whole row is removed and
# every later row shifts up one position, so 2025-10-20 becomes the new
# first day, the final day (2026-01-15) row disappears, and the sheet's
# used range shrinks from A1:C90 to A1:C89. Deleting the row also drops
# the now-unused "2025-10-19" shared string and renumbers every shared
# string after it, which is why the "Table" sheet's header cells (Issue /
# Validation / Pages) point at indices one lower than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete oldest-date row; Excel shifts the remaining rows up
# and recalculates the sheet's used range automatically.
$ws.Rows.Item(2).Delete()
